$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EG")

# Update modelled value for "Gas" (row 3)
$ws.Range("B3").Value = 145812.2776
$ws.Range("D3").Value = -15171.72240000003
$ws.Range("E3").Value = -54.6182006400001

# Update modelled value for "Solar PV" (row 5)
$ws.Range("B5").Value = 166.668
$ws.Range("D5").Value = -4339.332
$ws.Range("E5").Value = -15.6215952
